# Avatar position and path updates when a maze is regenerated.
# Avatar position calculated by physicsengine.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Color constants (BGR values as used by Excel COM Interior.Color)
$colorRed        = 255       # FFFF0000 - TODO / BROKEN (style s=1)
$colorGreen      = 5287936   # FF00B050 - DONE (style s=13)
$colorLightGreen = 5296274   # FF92D050 - SATISFACTORY / UNDERWAY (style s=12)

# Row 24: "Reset avatar position" is now DONE
$m24 = $ws.Range("M24")
$m24.Value = "DONE"
$m24.Interior.Color = $colorGreen

# Row 25: "Server::Instance()->Update() IS BROKEN" is now DONE
$m25 = $ws.Range("M25")
$m25.Value = "DONE"
$m25.Interior.Color = $colorGreen

# Row 26 (new): UsePhysics bool entry, status TODO
$ws.Range("K26").Value = "UsePhysics bool"
$ws.Range("L26").Value = "Allow client to change whether their position is calculated in the physics engine or not - include in ConnecetedClient struct"
$m26 = $ws.Range("M26")
$m26.Value = "TODO"
$m26.Interior.Color = $colorRed

# Row 21: add a note "R BUTTON ??" and mark status as SATISFACTORY
$ws.Range("L21").Value = "R BUTTON ??"
$m21 = $ws.Range("M21")
$m21.Value = "SATISFACTORY"
$m21.Interior.Color = $colorLightGreen

# Update the selected cell to reflect the new active cell
$ws.Range("M24").Select()
